$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append a new row (24) of mail-log data, mirroring the layout of the
# preceding rows (columns C and E are intentionally left blank).
$row = 24
$ws.Cells.Item($row, 1).Value = "Vraag over product"
$ws.Cells.Item($row, 2).Value = "documentatie@testbedrijf123.nl"
$ws.Cells.Item($row, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($row, 6).Value = "2025-08-19 21:07:36"
$ws.Cells.Item($row, 7).Value = "Nee"
$ws.Cells.Item($row, 8).Value = "Ja"
$ws.Cells.Item($row, 9).Value = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# The conditional-formatting rules on columns D, G, H, I and J were
# previously scoped to rows 2-23; extend them to cover the new row 24.
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range("$col`2:$col`23")
    $newRange = $ws.Range("$col`2:$col`24")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard sheet's aggregate count to reflect the new row.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 23
